$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H32").Value = 5297.8
$ws_ALC.Range("I32").Value = 3000
$ws_ALC.Range("J32").Value = 5872.25
$ws_ALC.Range("K32").Value = 3000
$ws_ALC.Range("L32").Value = 5872.25
$ws_ALC.Range("M32").Value = -2674
$ws_ALC.Range("N32").Value = -6524.25

$ws_ALC.Range("H43").Value = 1699.7142
$ws_ALC.Range("I43").Value = 1799.8
$ws_ALC.Range("K43").Value = 1799.8
$ws_ALC.Range("M43").Value = -1730.8

$ws_ALC.Range("H92").Value = 56279.39
$ws_ALC.Range("I92").Value = 77502.234
$ws_ALC.Range("K92").Value = 77502.234
$ws_ALC.Range("M92").Value = -76254.234

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H97").Value = 1594.1333
$ws_ARM.Range("I97").Value = 993.7143
$ws_ARM.Range("J97").Value = 10000
$ws_ARM.Range("K97").Value = 993.7143
$ws_ARM.Range("L97").Value = 10000
$ws_ARM.Range("M97").Value = -497.7143
$ws_ARM.Range("N97").Value = -10992

$ws_ARM.Range("H122").Value = 4000
$ws_ARM.Range("I122").Value = 4000
$ws_ARM.Range("K122").Value = 12000
$ws_ARM.Range("M122").Value = -9550

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 2397.7778
$ws_BSM.Range("I20").Value = 1947.5
$ws_BSM.Range("K20").Value = 1947.5
$ws_BSM.Range("M20").Value = -1700.5

$ws_BSM.Range("H134").Value = 8585.700000000001
$ws_BSM.Range("I134").Value = 8753.177
$ws_BSM.Range("K134").Value = 26259.531
$ws_BSM.Range("M134").Value = -23724.531

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H43").Value = 17475
$ws_CRP.Range("J43").Value = 17475
$ws_CRP.Range("L43").Value = 17475
$ws_CRP.Range("N43").Value = -17843

$ws_CRP.Range("H51").Value = 27499.5
$ws_CRP.Range("I51").Value = 27499.5
$ws_CRP.Range("K51").Value = 27499.5
$ws_CRP.Range("M51").Value = -26763.5

$ws_CRP.Range("H58").Value = 2024.0385
$ws_CRP.Range("I58").Value = 1596.5714
$ws_CRP.Range("K58").Value = 1596.5714
$ws_CRP.Range("M58").Value = -1393.5714

$ws_CRP.Range("H61").Value = 27499.5
$ws_CRP.Range("I61").Value = 27499.5
$ws_CRP.Range("K61").Value = 27499.5
$ws_CRP.Range("M61").Value = -27151.5

$ws_CRP.Range("H101").Value = 17475
$ws_CRP.Range("J101").Value = 17475
$ws_CRP.Range("L101").Value = 17475
$ws_CRP.Range("N101").Value = -23965

$ws_CRP.Range("H136").Value = 2024.0385
$ws_CRP.Range("I136").Value = 1596.5714
$ws_CRP.Range("K136").Value = 4789.7142
$ws_CRP.Range("M136").Value = -2239.7142

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H116").Value = 1500
$ws_CUL.Range("I116").Value = 1500
$ws_CUL.Range("K116").Value = 4500
$ws_CUL.Range("M116").Value = -1058

$ws_CUL.Range("H119").Value = 1665.3334
$ws_CUL.Range("I119").Value = 1665.3334
$ws_CUL.Range("K119").Value = 4996.0002
$ws_CUL.Range("M119").Value = -158.0002000000004

$ws_CUL.Range("H120").Value = 0
$ws_CUL.Range("I120").Value = 0
$ws_CUL.Range("K120").Value = 0
$ws_CUL.Range("M120").ClearContents()

$ws_CUL.Range("H121").Value = 716407.0600000001
$ws_CUL.Range("I121").Value = 839.9
$ws_CUL.Range("J121").Value = 2505325
$ws_CUL.Range("K121").Value = 2519.7
$ws_CUL.Range("L121").Value = 7515975
$ws_CUL.Range("M121").Value = -1209.7
$ws_CUL.Range("N121").Value = -7518595

$ws_CUL.Range("H131").Value = 3486.6667
$ws_CUL.Range("I131").Value = 3386
$ws_CUL.Range("J131").Value = 3990
$ws_CUL.Range("K131").Value = 10158
$ws_CUL.Range("L131").Value = 11970
$ws_CUL.Range("M131").Value = -5118
$ws_CUL.Range("N131").Value = -22050

$ws_CUL.Range("H133").Value = 19699.5
$ws_CUL.Range("I133").Value = 18497.5
$ws_CUL.Range("J133").Value = 20000
$ws_CUL.Range("K133").Value = 55492.5
$ws_CUL.Range("L133").Value = 60000
$ws_CUL.Range("M133").Value = -50432.5
$ws_CUL.Range("N133").Value = -70120

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H22").Value = 9500
$ws_GSM.Range("J22").Value = 9500
$ws_GSM.Range("L22").Value = 9500
$ws_GSM.Range("N22").Value = -10558

$ws_GSM.Range("H70").Value = 10301
$ws_GSM.Range("I70").Value = 9876.25
$ws_GSM.Range("J70").Value = 12000
$ws_GSM.Range("K70").Value = 9876.25
$ws_GSM.Range("L70").Value = 12000
$ws_GSM.Range("M70").Value = -9606.25
$ws_GSM.Range("N70").Value = -12540

$ws_GSM.Range("H73").Value = 10301
$ws_GSM.Range("I73").Value = 9876.25
$ws_GSM.Range("J73").Value = 12000
$ws_GSM.Range("K73").Value = 9876.25
$ws_GSM.Range("L73").Value = 12000
$ws_GSM.Range("M73").Value = -8940.25
$ws_GSM.Range("N73").Value = -13872

$ws_GSM.Range("H80").Value = 2664.125
$ws_GSM.Range("I80").Value = 2664.125
$ws_GSM.Range("J80").Value = 0
$ws_GSM.Range("K80").Value = 2664.125
$ws_GSM.Range("L80").Value = 0
$ws_GSM.Range("M80").ClearContents()
$ws_GSM.Range("N80").Value = -1666.125

$ws_GSM.Range("H83").Value = 2664.125
$ws_GSM.Range("I83").Value = 2664.125
$ws_GSM.Range("J83").Value = 0
$ws_GSM.Range("K83").Value = 13320.625
$ws_GSM.Range("L83").Value = 0
$ws_GSM.Range("M83").ClearContents()
$ws_GSM.Range("N83").Value = -8328.625

$ws_GSM.Range("H93").Value = 70251
$ws_GSM.Range("J93").Value = 70251
$ws_GSM.Range("L93").Value = 70251
$ws_GSM.Range("N93").Value = -73995

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H2").Value = 27499.25
$ws_LTW.Range("I2").Value = 20000
$ws_LTW.Range("K2").Value = 20000
$ws_LTW.Range("M2").Value = -19888

$ws_LTW.Range("H22").Value = 908.6667
$ws_LTW.Range("I22").Value = 863
$ws_LTW.Range("J22").Value = 1000
$ws_LTW.Range("K22").Value = 863
$ws_LTW.Range("L22").Value = 1000
$ws_LTW.Range("M22").Value = -568
$ws_LTW.Range("N22").Value = -1590

$ws_LTW.Range("H27").Value = 908.6667
$ws_LTW.Range("I27").Value = 863
$ws_LTW.Range("J27").Value = 1000
$ws_LTW.Range("K27").Value = 863
$ws_LTW.Range("L27").Value = 1000
$ws_LTW.Range("M27").Value = -756
$ws_LTW.Range("N27").Value = -1214

$ws_LTW.Range("H32").Value = 3500
$ws_LTW.Range("I32").Value = 2000
$ws_LTW.Range("J32").Value = 5000
$ws_LTW.Range("K32").Value = 2000
$ws_LTW.Range("L32").Value = 5000
$ws_LTW.Range("M32").Value = -1683
$ws_LTW.Range("N32").Value = -5634

$ws_LTW.Range("H46").Value = 3094.75
$ws_LTW.Range("I46").Value = 2560.125
$ws_LTW.Range("J46").Value = 3451.1667
$ws_LTW.Range("K46").Value = 2560.125
$ws_LTW.Range("L46").Value = 3451.1667
$ws_LTW.Range("M46").Value = -2372.125
$ws_LTW.Range("N46").Value = -3827.1667

$ws_LTW.Range("H61").Value = 2122.75
$ws_LTW.Range("J61").Value = 5150
$ws_LTW.Range("L61").Value = 5150
$ws_LTW.Range("N61").Value = -5554

$ws_LTW.Range("H93").Value = 1333.3914
$ws_LTW.Range("I93").Value = 1467.125
$ws_LTW.Range("J93").Value = 1027.7142
$ws_LTW.Range("K93").Value = 1467.125
$ws_LTW.Range("L93").Value = 1027.7142
$ws_LTW.Range("M93").Value = -219.125
$ws_LTW.Range("N93").Value = -3523.7142

$ws_LTW.Range("H113").Value = 2122.75
$ws_LTW.Range("J113").Value = 5150
$ws_LTW.Range("L113").Value = 5150
$ws_LTW.Range("N113").Value = -9490

$ws_LTW.Range("H136").Value = 10207626
$ws_LTW.Range("I136").Value = 2945.6047
$ws_LTW.Range("K136").Value = 8836.8141
$ws_LTW.Range("M136").Value = -6286.8141

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H62").Value = 5999.75
$ws_WVR.Range("I62").Value = 7332
$ws_WVR.Range("K62").Value = 7332
$ws_WVR.Range("M62").Value = -6708

$ws_WVR.Range("H65").Value = 5999.75
$ws_WVR.Range("I65").Value = 7332
$ws_WVR.Range("K65").Value = 36660
$ws_WVR.Range("M65").Value = -33540

$ws_WVR.Range("H81").Value = 8499.25
$ws_WVR.Range("I81").Value = 8499.25
$ws_WVR.Range("K81").Value = 16998.5
$ws_WVR.Range("M81").Value = -15937.5

$ws_WVR.Range("H84").Value = 8499.25
$ws_WVR.Range("I84").Value = 8499.25
$ws_WVR.Range("K84").Value = 84992.5
$ws_WVR.Range("M84").Value = -79688.5

$ws_WVR.Range("H132").Value = 1807.8
$ws_WVR.Range("I132").Value = 1473.2273
$ws_WVR.Range("J132").Value = 4261.3335
$ws_WVR.Range("K132").Value = 4419.6819
$ws_WVR.Range("L132").Value = 12784.0005
$ws_WVR.Range("M132").Value = -1889.6819
$ws_WVR.Range("N132").Value = -17844.0005

$ws_WVR.Range("H136").Value = 3891.2554
$ws_WVR.Range("I136").Value = 3122.4102
$ws_WVR.Range("J136").Value = 7639.375
$ws_WVR.Range("K136").Value = 9367.230599999999
$ws_WVR.Range("L136").Value = 22918.125
$ws_WVR.Range("M136").Value = -6817.230599999999
$ws_WVR.Range("N136").Value = -28018.125
